$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Feature / Importance data (rows 2-10), reordered and updated per target
$features = @(
    @("RSI", 0.4073932700996711),
    @("MACD", 0.2733529072579108),
    @("Signal_line", 0.1089292803744728),
    @("close_short", 0.04160514548460741),
    @("close_long", 0.04143784329497763),
    @("VIX_short", 0.03985920536463656),
    @("VIX", 0.03310259313622776),
    @("VIX_long", 0.02745191596034834),
    @("DJI", 0.02686783902714749)
)

for ($i = 0; $i -lt $features.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $features[$i][0]
    $ws.Cells.Item($row, 2).Value = $features[$i][1]
}
